$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2; existing rows 2-8 shift down to 3-9.
$ws.Rows.Item(2).Insert()

# The freshly inserted row picks up formatting from the row it was
# inserted at (bold/bordered header-ish style) - clear that so the new
# row matches the plain data-row look used by the rest of the table.
$ws.Range("A2:T2").ClearFormats()

# Column D holds dates; restore the same date number format used by the
# other rows in that column.
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Populate the new row 2 with its data.
$ws.Range("A2").Value = 11
$ws.Range("B2").Value = "Vega Monumental Concepción"
$ws.Range("C2").Value = "Bíobío"
$ws.Range("D2").Value = 44819
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100107
$ws.Range("H2").Value = "Otros"
$ws.Range("I2").Value = 100107011
$ws.Range("J2").Value = "Tuna"
$ws.Range("K2").Value = "Sin especificar"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 25000
$ws.Range("O2").Value = 26000
$ws.Range("P2").Value = 25500
$ws.Range("Q2").Value = "`$/caja 18 kilos granel"
$ws.Range("R2").Value = "Provincia de Melipilla"
$ws.Range("S2").Value = 1417
$ws.Range("T2").Value = 18
